$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 627.381
$ws.Range("I55").Value = 52.545456
$ws.Range("J55").Value = 1259.7
$ws.Range("K55").Value = 52.545456
$ws.Range("L55").Value = 1259.7
$ws.Range("M55").Value = 161.454544
$ws.Range("N55").Value = -1687.7
$ws.Range("H107").Value = 1679.5151
$ws.Range("J107").Value = 2229.5881
$ws.Range("L107").Value = 2229.5881
$ws.Range("N107").Value = -6069.5881

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6651.128
$ws.Range("I32").Value = 6296.2715
$ws.Range("J32").Value = 12399.8
$ws.Range("K32").Value = 6296.2715
$ws.Range("L32").Value = 12399.8
$ws.Range("M32").Value = -6009.2715
$ws.Range("N32").Value = -12973.8
$ws.Range("H45").Value = 3105.85
$ws.Range("I45").Value = 3185.923
$ws.Range("J45").Value = 2957.1428
$ws.Range("K45").Value = 3185.923
$ws.Range("L45").Value = 2957.1428
$ws.Range("M45").Value = -2808.923
$ws.Range("N45").Value = -3711.1428
$ws.Range("H61").Value = 9242.927
$ws.Range("I61").Value = 6902.385
$ws.Range("K61").Value = 6902.385
$ws.Range("M61").Value = -6690.385
$ws.Range("H63").Value = 2339
$ws.Range("I63").Value = 2355.5557
$ws.Range("J63").Value = 2190
$ws.Range("K63").Value = 2355.5557
$ws.Range("L63").Value = 2190
$ws.Range("M63").Value = -1669.5557
$ws.Range("N63").Value = -3562
$ws.Range("H66").Value = 2339
$ws.Range("I66").Value = 2355.5557
$ws.Range("J66").Value = 2190
$ws.Range("K66").Value = 11777.7785
$ws.Range("L66").Value = 10950
$ws.Range("M66").Value = -8345.7785
$ws.Range("N66").Value = -17814
$ws.Range("H74").Value = 3351.2068
$ws.Range("I74").Value = 1094.1082
$ws.Range("K74").Value = 1094.1082
$ws.Range("M74").Value = -220.1081999999999
$ws.Range("H77").Value = 3351.2068
$ws.Range("I77").Value = 1094.1082
$ws.Range("K77").Value = 5470.540999999999
$ws.Range("M77").Value = -1102.540999999999
$ws.Range("H110").Value = 807.65
$ws.Range("I110").Value = 686.2778
$ws.Range("K110").Value = 686.2778
$ws.Range("M110").Value = 1358.7222
$ws.Range("H122").Value = 3324.8845
$ws.Range("I122").Value = 2659.762
$ws.Range("J122").Value = 6118.4
$ws.Range("K122").Value = 7979.286
$ws.Range("L122").Value = 18355.2
$ws.Range("M122").Value = -5529.286
$ws.Range("N122").Value = -23255.2
$ws.Range("H132").Value = 1648.386
$ws.Range("I132").Value = 1566.3024
$ws.Range("K132").Value = 4698.9072
$ws.Range("M132").Value = -2168.9072
$ws.Range("H136").Value = 9242.927
$ws.Range("I136").Value = 6902.385
$ws.Range("K136").Value = 20707.155
$ws.Range("M136").Value = -18157.155

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2655.182
$ws.Range("I99").Value = 1889.8334
$ws.Range("J99").Value = 6099.25
$ws.Range("K99").Value = 1889.8334
$ws.Range("L99").Value = 6099.25
$ws.Range("M99").Value = -391.8334
$ws.Range("N99").Value = -9095.25
$ws.Range("H105").Value = 3162.718
$ws.Range("I105").Value = 2615.6287
$ws.Range("K105").Value = 2615.6287
$ws.Range("M105").Value = -868.6287000000002
$ws.Range("H107").Value = 863.05
$ws.Range("I107").Value = 680.1177
$ws.Range("J107").Value = 1899.6666
$ws.Range("K107").Value = 680.1177
$ws.Range("L107").Value = 1899.6666
$ws.Range("M107").Value = 1239.8823
$ws.Range("N107").Value = -5739.6666
$ws.Range("H132").Value = 100000
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 100000
$ws.Range("N132").Value = -110120
$ws.Range("H134").Value = 6860.533
$ws.Range("I134").Value = 3040.6667
$ws.Range("K134").Value = 9122.000100000001
$ws.Range("M134").Value = -6587.000100000001
$ws.Range("H140").Value = 175148
$ws.Range("J140").Value = 194177.6
$ws.Range("L140").Value = 194177.6
$ws.Range("N140").Value = -204537.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3264.639
$ws.Range("I31").Value = 2228
$ws.Range("K31").Value = 2228
$ws.Range("M31").Value = -1933
$ws.Range("H34").Value = 3264.639
$ws.Range("I34").Value = 2228
$ws.Range("K34").Value = 2228
$ws.Range("M34").Value = -2026
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H58").Value = 2821.0645
$ws.Range("I58").Value = 2358.125
$ws.Range("J58").Value = 3314.8667
$ws.Range("K58").Value = 2358.125
$ws.Range("L58").Value = 3314.8667
$ws.Range("M58").Value = -2155.125
$ws.Range("N58").Value = -3720.8667
$ws.Range("H105").Value = 1984.375
$ws.Range("I105").Value = 728.4
$ws.Range("J105").Value = 4077.6667
$ws.Range("K105").Value = 728.4
$ws.Range("L105").Value = 4077.6667
$ws.Range("M105").Value = 1018.6
$ws.Range("N105").Value = -7571.6667
$ws.Range("H132").Value = 955028.9399999999
$ws.Range("I132").Value = 1601839
$ws.Range("K132").Value = 4805517
$ws.Range("M132").Value = -4802987
$ws.Range("H134").Value = 5459.0557
$ws.Range("I134").Value = 2251.9285
$ws.Range("K134").Value = 6755.7855
$ws.Range("M134").Value = -4220.7855
$ws.Range("H136").Value = 2821.0645
$ws.Range("I136").Value = 2358.125
$ws.Range("J136").Value = 3314.8667
$ws.Range("K136").Value = 7074.375
$ws.Range("L136").Value = 9944.6001
$ws.Range("M136").Value = -4524.375
$ws.Range("N136").Value = -15044.6001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1760.96
$ws.Range("I12").Value = 2835.7
$ws.Range("J12").Value = 1044.4667
$ws.Range("K12").Value = 8507.099999999999
$ws.Range("L12").Value = 3133.4001
$ws.Range("M12").Value = -8334.099999999999
$ws.Range("N12").Value = -3479.4001
$ws.Range("H36").Value = 2183.6667
$ws.Range("I36").Value = 1412.625
$ws.Range("J36").Value = 3725.75
$ws.Range("K36").Value = 4237.875
$ws.Range("L36").Value = 11177.25
$ws.Range("M36").Value = -4068.875
$ws.Range("N36").Value = -11515.25
$ws.Range("H127").Value = 2698.5
$ws.Range("J127").Value = 2698.5
$ws.Range("L127").Value = 8095.5
$ws.Range("N127").Value = -18015.5
$ws.Range("H132").Value = 3019.0278
$ws.Range("I132").Value = 1496.8572
$ws.Range("J132").Value = 3386.4482
$ws.Range("K132").Value = 13471.7148
$ws.Range("L132").Value = 30478.0338
$ws.Range("M132").Value = -10941.7148
$ws.Range("N132").Value = -35538.0338
$ws.Range("H134").Value = 2580.75
$ws.Range("I134").Value = 2291.6924
$ws.Range("K134").Value = 6875.0772
$ws.Range("M134").Value = -1805.0772
$ws.Range("H138").Value = 5794536.5
$ws.Range("I138").Value = 1255899.2
$ws.Range("K138").Value = 3767697.6
$ws.Range("M138").Value = -3762557.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 29617.334
$ws.Range("J48").Value = 29617.334
$ws.Range("L48").Value = 29617.334
$ws.Range("N48").Value = -30587.334
$ws.Range("H123").Value = 45957
$ws.Range("J123").Value = 45957
$ws.Range("L123").Value = 45957
$ws.Range("N123").Value = -50857
$ws.Range("H132").Value = 2437.5476
$ws.Range("I132").Value = 2298.1082
$ws.Range("K132").Value = 6894.3246
$ws.Range("M132").Value = -4364.3246

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1847.6552
$ws.Range("I22").Value = 561.2222
$ws.Range("K22").Value = 561.2222
$ws.Range("M22").Value = -266.2222
$ws.Range("H27").Value = 1847.6552
$ws.Range("I27").Value = 561.2222
$ws.Range("K27").Value = 561.2222
$ws.Range("M27").Value = -454.2222
$ws.Range("H46").Value = 10591.833
$ws.Range("I46").Value = 4910
$ws.Range("J46").Value = 12777.154
$ws.Range("K46").Value = 4910
$ws.Range("L46").Value = 12777.154
$ws.Range("M46").Value = -4722
$ws.Range("N46").Value = -13153.154
$ws.Range("H122").Value = 4693.1055
$ws.Range("I122").Value = 4166.4614
$ws.Range("K122").Value = 12499.3842
$ws.Range("M122").Value = -10049.3842
$ws.Range("H124").Value = 49500
$ws.Range("J124").Value = 49500
$ws.Range("L124").Value = 49500
$ws.Range("N124").Value = -59320

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 26074
$ws.Range("J47").Value = 26074
$ws.Range("L47").Value = 26074
$ws.Range("N47").Value = -27218
$ws.Range("H51").Value = 41999
$ws.Range("I51").Value = 37998.5
$ws.Range("J51").Value = 50000
$ws.Range("K51").Value = 37998.5
$ws.Range("L51").Value = 50000
$ws.Range("M51").Value = -37488.5
$ws.Range("N51").Value = -51020
$ws.Range("H52").Value = 51150
$ws.Range("J52").Value = 52500
$ws.Range("L52").Value = 52500
$ws.Range("N52").Value = -52952
$ws.Range("H54").Value = 12000
$ws.Range("I54").Value = 3000
$ws.Range("J54").Value = 21000
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 21000
$ws.Range("M54").Value = -2480
$ws.Range("N54").Value = -22040
$ws.Range("H100").Value = 1191.0741
$ws.Range("I100").Value = 1054.2273
$ws.Range("K100").Value = 2108.4546
$ws.Range("M100").Value = -1567.4546
$ws.Range("H122").Value = 3060.8462
$ws.Range("I122").Value = 3379.8
$ws.Range("J122").Value = 1997.6666
$ws.Range("K122").Value = 10139.4
$ws.Range("L122").Value = 5992.9998
$ws.Range("M122").Value = -7689.400000000001
$ws.Range("N122").Value = -10892.9998
$ws.Range("H126").Value = 2219.2856
$ws.Range("I126").Value = 2139.5908
$ws.Range("K126").Value = 6418.7724
$ws.Range("M126").Value = -3948.7724
$ws.Range("H136").Value = 4537.887
$ws.Range("I136").Value = 4662.849
$ws.Range("K136").Value = 13988.547
$ws.Range("M136").Value = -11438.547
